$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.04155
$ws.Range("H2").Value = 0.12465
$ws.Range("I2").Value = 0.0001466168179836329
$ws.Range("J2").Value = 0.0001466168179836329
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.007501333333333333
$ws.Range("N2").Value = 0.022504
$ws.Range("O2").Value = 0.1758070060310615
$ws.Range("P2").Value = 0.1758070060310615
$ws.Range("Q2").Value = 0.0003116804
$ws.Range("R2").Value = 0.0028051236
$ws.Range("S2").Value = 0.00002577626380350359
$ws.Range("T2").Value = 0.00002577626380350359

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.04155
$ws.Range("H3").Value = 0.12465
$ws.Range("I3").Value = 0.0001466168179836329
$ws.Range("J3").Value = 0.0001466168179836329
$ws.Range("O3").Value = 0.8241929939689385
$ws.Range("P3").Value = 0.8241929939689384
$ws.Range("Q3").Value = 0.001461175
$ws.Range("R3").Value = 0.013150575
$ws.Range("S3").Value = 0.0001208405541801293
$ws.Range("T3").Value = 0.0001208405541801293

# Row 4
$ws.Range("I4").Value = 0.9992428949822291
$ws.Range("J4").Value = 0.9992428949822291
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.007501333333333333
$ws.Range("N4").Value = 0.022504
$ws.Range("O4").Value = 0.1758070060310615
$ws.Range("P4").Value = 0.1758070060310615
$ws.Range("Q4").Value = 2.124206687120889
$ws.Range("R4").Value = 19.117860184088
$ws.Range("S4").Value = 0.1756739016646361
$ws.Range("T4").Value = 0.1756739016646361

# Row 5
$ws.Range("I5").Value = 0.9992428949822291
$ws.Range("J5").Value = 0.9992428949822291
$ws.Range("O5").Value = 0.8241929939689385
$ws.Range("P5").Value = 0.8241929939689384
$ws.Range("S5").Value = 0.823568993317593
$ws.Range("T5").Value = 0.8235689933175929

# Row 6
$ws.Range("I6").Value = 0.0006104881997874136
$ws.Range("J6").Value = 0.0006104881997874135
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.007501333333333333
$ws.Range("N6").Value = 0.022504
$ws.Range("O6").Value = 0.1758070060310615
$ws.Range("P6").Value = 0.1758070060310615
$ws.Range("Q6").Value = 0.001297785676444444
$ws.Range("R6").Value = 0.011680071088
$ws.Range("S6").Value = 0.0001073281026219177
$ws.Range("T6").Value = 0.0001073281026219177

# Row 7
$ws.Range("I7").Value = 0.0006104881997874136
$ws.Range("J7").Value = 0.0006104881997874135
$ws.Range("O7").Value = 0.8241929939689385
$ws.Range("P7").Value = 0.8241929939689384
$ws.Range("S7").Value = 0.0005031600971654959
$ws.Range("T7").Value = 0.0005031600971654957
